$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.890.45'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '2.116.56'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '348.12'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5181'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4465'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.91'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09350'
$ws.Range('E10').Value = '  +4.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.183'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.24'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '2.099.02'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.841'
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.315'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '102.68'
$ws.Range('E16').Value = '  +3.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001165'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.007'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.52'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06671'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.309'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').Value = '29.930.99'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.327'
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.361.49'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.19'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.556'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.49'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.02'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.151'
$ws.Range('E31').Value = '  -2.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.792'
$ws.Range('E32').Value = '  +8.77%  '
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.253'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.967'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.394'
$ws.Range('E36').Value = '  +7.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.84'
$ws.Range('E37').Value = '  +6.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02595'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06807'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.74'
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7025'
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.344'
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2243'
$ws.Range('E43').Value = '  -2.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6862'
$ws.Range('E44').Value = '  +7.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.50'
$ws.Range('E45').Value = '  +1.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.364'
$ws.Range('E46').Value = '  +3.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.005'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000353'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.219'
$ws.Range('E50').Value = '  +4.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.226'
$ws.Range('E51').Value = '  +0.64%  '
